# Workbook/settings refresh + data cleanup on Sheet1.
# - Fixes the "AI In.Imit_Mode" tag string: the six duplicate TAG entries
#   (B6, B8, B10, B12, B14, B16) that were wrongly carrying the shared
#   random-suffix formula are cleared out (these rows have no real tag name,
#   hence "No Tag Name" — they should stay blank like their B3/B5/B7/... siblings).
# - Moves the saved selection to C18 (and off the scrolled A153 view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("B16").ClearContents()

$ws.Activate()
$ws.Range("C18").Select()
